$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "InstantPay " / "(instant transactions)..." bullet -> merged, reworded
#    bullet that now shares the numId=2 bulleted list, with updated run
#    formatting (Open Sans / color 252525) and single-line spacing.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("InstantPay*confirm a payment.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the InstantPay bullet paragraph"
}
$para1 = $rng1.Paragraphs(1)

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr>' +
          '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
          '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
          '<w:textAlignment w:val="baseline"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="252525"/>' +
            '<w:sz w:val="21"/>' +
            '<w:szCs w:val="21"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
            '<w:color w:val="252525"/>' +
            '<w:sz w:val="21"/>' +
            '<w:szCs w:val="21"/>' +
          '</w:rPr>' +
          '<w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t>' +
        '</w:r>' +
        '</w:p>'

$para1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Will have more services added later" bullet -> reworded SmartRewards
#    bullet, same numId=2 list, updated run formatting and single-line
#    spacing.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Will have more services added later", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the 'Will have more services' bullet paragraph"
}
$para2 = $rng2.Paragraphs(1)

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr>' +
          '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
          '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
          '<w:textAlignment w:val="baseline"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
            '<w:color w:val="252525"/>' +
            '<w:sz w:val="21"/>' +
            '<w:szCs w:val="21"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
            '<w:color w:val="252525"/>' +
            '<w:sz w:val="21"/>' +
            '<w:szCs w:val="21"/>' +
          '</w:rPr>' +
          '<w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t>' +
        '</w:r>' +
        '</w:p>'

$para2.Range.InsertXML($xml2)
